# Update scraped view/like counters across the four sheets of
# 北京-漫展信息.xlsx (gh-pages data refresh @ 456a3b4).
#
# Sheet tab order: 展览(1), 演出(2), 本地生活(3), 全部类型(4) -- the last
# sheet is a merged view of the first three, so several values are
# touched twice (once in their "home" sheet, once in the merged copy).

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetShow       = $wb.Worksheets.Item("演出")
$sheetLocal      = $wb.Worksheets.Item("本地生活")
$sheetAll        = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibition)
$sheetExhibition.Range("F4").Value  = 33
$sheetExhibition.Range("F7").Value  = 2171
$sheetExhibition.Range("F13").Value = 1481
$sheetExhibition.Range("F15").Value = 566
$sheetExhibition.Range("F16").Value = 413
$sheetExhibition.Range("F17").Value = 413
$sheetExhibition.Range("F18").Value = 789
$sheetExhibition.Range("F19").Value = 461
$sheetExhibition.Range("F20").Value = 2991
$sheetExhibition.Range("F22").Value = 117
$sheetExhibition.Range("F25").Value = 551
$sheetExhibition.Range("F26").Value = 251
$sheetExhibition.Range("F27").Value = 1002
$sheetExhibition.Range("F28").Value = 737
$sheetExhibition.Range("F30").Value = 777
$sheetExhibition.Range("F31").Value = 758

# 演出 (Show)
$sheetShow.Range("G3").Value  = 108
$sheetShow.Range("F20").Value = 203
$sheetShow.Range("F21").Value = 144
$sheetShow.Range("F22").Value = 450

# 本地生活 (Local life)
$sheetLocal.Range("F4").Value = 379

# 全部类型 (All types, merged copy of the three sheets above)
$sheetAll.Range("F6").Value  = 379
$sheetAll.Range("F8").Value  = 33
$sheetAll.Range("G9").Value  = 108
$sheetAll.Range("F14").Value = 2171
$sheetAll.Range("F24").Value = 1481
$sheetAll.Range("F25").Value = 1481
$sheetAll.Range("F28").Value = 566
$sheetAll.Range("F29").Value = 413
$sheetAll.Range("F30").Value = 413
$sheetAll.Range("F32").Value = 789
$sheetAll.Range("F33").Value = 461
$sheetAll.Range("F35").Value = 2991
$sheetAll.Range("F36").Value = 117
$sheetAll.Range("F40").Value = 551
$sheetAll.Range("F41").Value = 251
$sheetAll.Range("F42").Value = 1002
$sheetAll.Range("F45").Value = 203
$sheetAll.Range("F46").Value = 144
$sheetAll.Range("F47").Value = 450
$sheetAll.Range("F48").Value = 737
$sheetAll.Range("F50").Value = 777
$sheetAll.Range("F51").Value = 758
